# Scheduled-runner price/profit refresh: pushes newly-pulled
# Universalis market-board averages (H/I/J/K/L) and recomputed
# profit deltas (M/N) into each job sheets leve table.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 17783510
$ws.Range("J62").Value = 5449
$ws.Range("L62").Value = 5449
$ws.Range("N62").Value = -6697

$ws.Range("H65").Value = 17783510
$ws.Range("J65").Value = 5449
$ws.Range("L65").Value = 27245
$ws.Range("N65").Value = -33485

$ws.Range("H100").Value = 6768.1333
$ws.Range("I100").Value = 1099.375
$ws.Range("J100").Value = 7993.811
$ws.Range("K100").Value = 1099.375
$ws.Range("L100").Value = 7993.811
$ws.Range("M100").Value = -558.375
$ws.Range("N100").Value = -9075.811

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3367855
$ws.Range("I2").Value = 4329862
$ws.Range("K2").Value = 4329862
$ws.Range("M2").Value = -4329749

$ws.Range("H21").Value = 3465.6667
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").ClearContents()

$ws.Range("H74").Value = 7704
$ws.Range("I74").Value = 1649.8
$ws.Range("J74").Value = 13758.2
$ws.Range("K74").Value = 1649.8
$ws.Range("L74").Value = 13758.2
$ws.Range("M74").Value = -775.8
$ws.Range("N74").Value = -15506.2

$ws.Range("H77").Value = 7704
$ws.Range("I77").Value = 1649.8
$ws.Range("J77").Value = 13758.2
$ws.Range("K77").Value = 8249
$ws.Range("L77").Value = 68791
$ws.Range("M77").Value = -3881
$ws.Range("N77").Value = -77527

$ws.Range("H102").Value = 35715880
$ws.Range("I102").Value = 1859.8334
$ws.Range("K102").Value = 1859.8334
$ws.Range("M102").Value = -237.8334

$ws.Range("H110").Value = 12501947
$ws.Range("I110").Value = 27778440
$ws.Range("K110").Value = 27778440
$ws.Range("M110").Value = -27776395

$ws.Range("H116").Value = 3367855
$ws.Range("I116").Value = 4329862
$ws.Range("K116").Value = 4329862
$ws.Range("M116").Value = -4327568

$ws.Range("H122").Value = 2195
$ws.Range("I122").Value = 2369.4443
$ws.Range("K122").Value = 7108.3329
$ws.Range("M122").Value = -4658.3329

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3367855
$ws.Range("I3").Value = 4329862
$ws.Range("K3").Value = 4329862
$ws.Range("M3").Value = -4329748

$ws.Range("H99").Value = 1657
$ws.Range("I99").Value = 1339.8
$ws.Range("J99").Value = 2450
$ws.Range("K99").Value = 1339.8
$ws.Range("L99").Value = 2450
$ws.Range("M99").Value = 158.2
$ws.Range("N99").Value = -5446

$ws.Range("H140").Value = 98665.664
$ws.Range("J140").Value = 98665.664
$ws.Range("L140").Value = 98665.664
$ws.Range("N140").Value = -109025.664

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4030.8538
$ws.Range("J31").Value = 6129.7896
$ws.Range("L31").Value = 6129.7896
$ws.Range("N31").Value = -6719.7896

$ws.Range("H34").Value = 4030.8538
$ws.Range("J34").Value = 6129.7896
$ws.Range("L34").Value = 6129.7896
$ws.Range("N34").Value = -6533.7896

$ws.Range("H59").Value = 61267.5
$ws.Range("J59").Value = 102449.5
$ws.Range("L59").Value = 102449.5
$ws.Range("N59").Value = -104739.5

$ws.Range("H68").Value = 171142.86
$ws.Range("J68").Value = 79600
$ws.Range("L68").Value = 79600
$ws.Range("N68").Value = -81098

$ws.Range("H71").Value = 171142.86
$ws.Range("J71").Value = 79600
$ws.Range("L71").Value = 238800
$ws.Range("N71").Value = -246288

$ws.Range("H105").Value = 1794.5
$ws.Range("I105").Value = 1403.762
$ws.Range("K105").Value = 1403.762
$ws.Range("M105").Value = 343.2380000000001

$ws.Range("H107").Value = 125000610
$ws.Range("I107").Value = 125000610
$ws.Range("K107").Value = 125000610
$ws.Range("M107").Value = -124998690

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H62").Value = 1268.119
$ws.Range("J62").Value = 3500
$ws.Range("L62").Value = 10500
$ws.Range("N62").Value = -11872

$ws.Range("H65").Value = 1268.119
$ws.Range("J65").Value = 3500
$ws.Range("L65").Value = 31500
$ws.Range("N65").Value = -38364

$ws.Range("H109").Value = 3131.4285
$ws.Range("I109").Value = 3131.4285
$ws.Range("K109").Value = 9394.2855
$ws.Range("M109").Value = -8354.2855

$ws.Range("H129").Value = 35715504
$ws.Range("I129").Value = 544.8570999999999
$ws.Range("J129").Value = 71430460
$ws.Range("K129").Value = 1634.5713
$ws.Range("L129").Value = 214291380
$ws.Range("M129").Value = 3365.4287
$ws.Range("N129").Value = -214301380

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H6").Value = 5031
$ws.Range("I6").Value = 1500
$ws.Range("J6").Value = 6796.5
$ws.Range("K6").Value = 1500
$ws.Range("L6").Value = 6796.5
$ws.Range("M6").Value = -1387
$ws.Range("N6").Value = -7022.5

$ws.Range("H7").Value = 5000750
$ws.Range("J7").Value = 5000750
$ws.Range("L7").Value = 5000750
$ws.Range("N7").Value = -5000974

$ws.Range("H8").Value = 5000750
$ws.Range("J8").Value = 5000750
$ws.Range("L8").Value = 5000750
$ws.Range("N8").Value = -5001028

$ws.Range("H14").Value = 428587.1
$ws.Range("J14").Value = 14888.333
$ws.Range("L14").Value = 14888.333
$ws.Range("N14").Value = -15224.333

$ws.Range("H16").Value = 5031
$ws.Range("I16").Value = 1500
$ws.Range("J16").Value = 6796.5
$ws.Range("K16").Value = 1500
$ws.Range("L16").Value = 6796.5
$ws.Range("M16").Value = -1250
$ws.Range("N16").Value = -7296.5

$ws.Range("H19").Value = 18077.75
$ws.Range("I19").Value = 19999
$ws.Range("J19").Value = 17903.092
$ws.Range("K19").Value = 19999
$ws.Range("L19").Value = 17903.092
$ws.Range("M19").Value = -19711
$ws.Range("N19").Value = -18479.092

$ws.Range("H132").Value = 6598.96
$ws.Range("I132").Value = 4689.2
$ws.Range("J132").Value = 9463.6
$ws.Range("K132").Value = 14067.6
$ws.Range("L132").Value = 28390.8
$ws.Range("M132").Value = -11537.6
$ws.Range("N132").Value = -33450.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("N40").ClearContents()

$ws.Range("H61").Value = 2124.913
$ws.Range("I61").Value = 2429.3157
$ws.Range("J61").Value = 679
$ws.Range("K61").Value = 2429.3157
$ws.Range("L61").Value = 679
$ws.Range("M61").Value = -2227.3157
$ws.Range("N61").Value = -1083

$ws.Range("H100").Value = 35717550
$ws.Range("J100").Value = 3935.6
$ws.Range("L100").Value = 3935.6
$ws.Range("N100").Value = -5017.6

$ws.Range("H113").Value = 2124.913
$ws.Range("I113").Value = 2429.3157
$ws.Range("J113").Value = 679
$ws.Range("K113").Value = 2429.3157
$ws.Range("L113").Value = 679
$ws.Range("M113").Value = -259.3157000000001
$ws.Range("N113").Value = -5019

$ws.Range("H132").Value = 3883.5
$ws.Range("I132").Value = 3431.093
$ws.Range("J132").Value = 4729.304
$ws.Range("K132").Value = 10293.279
$ws.Range("L132").Value = 14187.912
$ws.Range("M132").Value = -7763.278999999999
$ws.Range("N132").Value = -19247.912

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H69").Value = 16266.4
$ws.Range("J69").Value = 16266.4
$ws.Range("L69").Value = 16266.4
$ws.Range("N69").Value = -17764.4

$ws.Range("H72").Value = 16266.4
$ws.Range("J72").Value = 16266.4
$ws.Range("L72").Value = 48799.2
$ws.Range("N72").Value = -56287.2

$ws.Range("H81").Value = 11931.405
$ws.Range("I81").Value = 6047.2856
$ws.Range("J81").Value = 15513.044
$ws.Range("K81").Value = 12094.5712
$ws.Range("L81").Value = 31026.088
$ws.Range("M81").Value = -11033.5712
$ws.Range("N81").Value = -33148.088

$ws.Range("H84").Value = 11931.405
$ws.Range("I84").Value = 6047.2856
$ws.Range("J84").Value = 15513.044
$ws.Range("K84").Value = 60472.856
$ws.Range("L84").Value = 155130.44
$ws.Range("M84").Value = -55168.856
$ws.Range("N84").Value = -165738.44

$ws.Range("H113").Value = 799.2174
$ws.Range("I113").Value = 1033
$ws.Range("J113").Value = 435.55554
$ws.Range("K113").Value = 3099
$ws.Range("L113").Value = 1306.66662
$ws.Range("M113").Value = -929
$ws.Range("N113").Value = -5646.66662

$ws.Range("H132").Value = 6285.6924
$ws.Range("I132").Value = 5745.9585
$ws.Range("J132").Value = 12762.5
$ws.Range("K132").Value = 17237.8755
$ws.Range("L132").Value = 38287.5
$ws.Range("M132").Value = -14707.8755
$ws.Range("N132").Value = -43347.5
